$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2 = 'Tee,Kaftan'
    3 = 'Kaftan,Tee'
    4 = 'Halter,Sweatpants'
    5 = 'Parka,Jumpsuit'
    6 = 'Trunks,Jodhpurs'
    7 = 'Jumpsuit,Dress'
    8 = 'Jumpsuit,Blouse'
    9 = 'Jumpsuit,Kaftan'
    10 = 'Cutoffs,Halter'
    11 = 'Caftan,Jumpsuit'
    12 = 'Jumpsuit,Kaftan'
    13 = 'Tee,Kaftan'
    14 = 'Jumpsuit,Blazer'
    15 = 'Jumpsuit,Caftan'
    16 = 'Jumpsuit,Blouse'
    17 = 'Kaftan,Parka'
    18 = 'Tee,Parka'
    19 = 'Blouse,Trunks'
    20 = 'Blouse,Parka'
    21 = 'Tee,Kaftan'
    22 = 'Blouse,Jumpsuit'
    23 = 'Tee,Kaftan'
    24 = 'Blouse,Trunks'
    25 = 'Sweatpants,Blouse'
    26 = 'Kaftan,Tee'
    27 = 'Jumpsuit,Halter'
    28 = 'Jumpsuit,Halter'
    29 = 'Jodhpurs,Parka'
    30 = 'Kaftan,Tee'
    31 = 'Trunks,Sweatpants'
    32 = 'Tee,Kaftan'
    33 = 'Tee,Kaftan'
    34 = 'Blazer,Halter'
    35 = 'Halter,Coat'
    36 = 'Jumpsuit,Halter'
    37 = 'Jumpsuit,Blouse'
    38 = 'Tee,Kaftan'
    39 = 'Blazer,Jumpsuit'
    40 = 'Parka,Blouse'
    41 = 'Jumpsuit,Halter'
    42 = 'Kaftan,Tee'
    43 = 'Tee,Kaftan'
    44 = 'Tee,Kaftan'
    45 = 'Tee,Kaftan'
    46 = 'Tee,Kaftan'
    47 = 'Tee,Kaftan'
    48 = 'Tee,Kaftan'
    49 = 'Tee,Parka'
    50 = 'Tee,Kaftan'
    51 = 'Blouse,Sweatpants'
    52 = 'Jumpsuit,Kaftan'
    53 = 'Parka,Jumpsuit'
    54 = 'Jumpsuit,Blouse'
    55 = 'Jumpsuit,Blouse'
    56 = 'Parka,Jodhpurs'
    57 = 'Trunks,Jodhpurs'
    58 = 'Tee,Kaftan'
    59 = 'Parka,Blouse'
    60 = 'Jumpsuit,Blouse'
    61 = 'Kaftan,Tee'
    62 = 'Jumpsuit,Halter'
    63 = 'Jumpsuit,Trunks'
    64 = 'Blouse,Jumpsuit'
    65 = 'Jumpsuit,Halter'
    66 = 'Blouse,Parka'
    67 = 'Jumpsuit,Dress'
    68 = 'Jumpsuit,Blazer'
    69 = 'Blazer,Halter'
    70 = 'Dress,Sweatpants'
    71 = 'Blouse,Parka'
    72 = 'Dress,Kaftan'
    73 = 'Tee,Kaftan'
    74 = 'Dress,Parka'
    75 = 'Trunks,Jodhpurs'
    76 = 'Jumpsuit,Kaftan'
    77 = 'Trunks,Caftan'
    78 = 'Parka,Trunks'
    79 = 'Parka,Top'
    80 = 'Blazer,Blouse'
    81 = 'Tee,Kaftan'
    82 = 'Jumpsuit,Kaftan'
    83 = 'Jumpsuit,Dress'
    84 = 'Jumpsuit,Blouse'
    85 = 'Jumpsuit,Halter'
    86 = 'Parka,Dress'
    87 = 'Jumpsuit,Kaftan'
    88 = 'Trunks,Dress'
    89 = 'Halter,Parka'
    90 = 'Jumpsuit,Blazer'
    91 = 'Halter,Jumpsuit'
    92 = 'Jodhpurs,Trunks'
    93 = 'Kaftan,Tee'
    94 = 'Blouse,Parka'
    95 = 'Tee,Jumpsuit'
    96 = 'Blouse,Jumpsuit'
    97 = 'Blouse,Parka'
    98 = 'Blouse,Dress'
    99 = 'Caftan,Trunks'
    100 = 'Kaftan,Tee'
    101 = 'Jumpsuit,Halter'
    102 = 'Dress,Kaftan'
    103 = 'Blazer,Top'
    104 = 'Kaftan,Tee'
    105 = 'Kaftan,Tee'
    106 = 'Kaftan,Tee'
    107 = 'Tee,Kaftan'
    108 = 'Trunks,Jodhpurs'
    109 = 'Blouse,Jumpsuit'
    110 = 'Tee,Kaftan'
    111 = 'Kaftan,Tee'
    112 = 'Kaftan,Tee'
    113 = 'Kaftan,Tee'
    114 = 'Parka,Caftan'
    115 = 'Jumpsuit,Tee'
    116 = 'Jumpsuit,Parka'
    117 = 'Jumpsuit,Tee'
    118 = 'Kaftan,Jodhpurs'
    119 = 'Trunks,Turtleneck'
    120 = 'Jumpsuit,Trunks'
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $values[$row]
}
